# Auto-generated Excel COM-interop script to apply market-data refresh
# to the Leve profit calc sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each changed row we update currentAveragePrice / NQ / HQ columns
# (H-N) to the freshly scraped values.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 18
$ws.Range("H18").Value = 964.5
$ws.Range("I18").Value = 416.1111
$ws.Range("K18").Value = 416.1111
$ws.Range("M18").Value = -132.1111

# Row 38
$ws.Range("H38").Value = 4828.6313
$ws.Range("J38").Value = 8569.714
$ws.Range("L38").Value = 25709.142
$ws.Range("N38").Value = -26453.142

# Row 40
$ws.Range("H40").Value = 6523
$ws.Range("J40").Value = 7394
$ws.Range("L40").Value = 7394
$ws.Range("N40").Value = -7744

# Row 103
$ws.Range("H103").Value = 912.82355
$ws.Range("J103").Value = 845.3125
$ws.Range("L103").Value = 2535.9375
$ws.Range("N103").Value = -3707.9375

# Row 115
$ws.Range("H115").Value = 393
$ws.Range("I115").Value = 393
$ws.Range("K115").Value = 1179
$ws.Range("M115").Value = 388

# Row 132
$ws.Range("H132").Value = 4951.593
$ws.Range("I132").Value = 5108.731
$ws.Range("J132").Value = 866
$ws.Range("K132").Value = 15326.193
$ws.Range("L132").Value = 2598
$ws.Range("M132").Value = -12796.193
$ws.Range("N132").Value = -7658

# Row 138
$ws.Range("H138").Value = 4687.424
$ws.Range("I138").Value = 1381.84
$ws.Range("K138").Value = 4145.52
$ws.Range("M138").Value = 994.4800000000005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# Row 23
$ws.Range("H23").Value = 6166.6665
$ws.Range("I23").Value = 8000
$ws.Range("J23").Value = 5250
$ws.Range("K23").Value = 8000
$ws.Range("L23").Value = 5250
$ws.Range("M23").Value = -7741
$ws.Range("N23").Value = -5768

# Row 32
$ws.Range("H32").Value = 3695.365
$ws.Range("I32").Value = 3163.386
$ws.Range("J32").Value = 8749.166999999999
$ws.Range("K32").Value = 3163.386
$ws.Range("L32").Value = 8749.166999999999
$ws.Range("M32").Value = -2876.386
$ws.Range("N32").Value = -9323.166999999999

# Row 37
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 22
$ws.Range("H22").Value = 337.5
$ws.Range("I22").Value = 350
$ws.Range("K22").Value = 350
$ws.Range("M22").Value = -177

# Row 86
$ws.Range("H86").Value = 18787.715
$ws.Range("I86").Value = 23303
$ws.Range("J86").Value = 7499.5
$ws.Range("K86").Value = 23303
$ws.Range("L86").Value = 7499.5
$ws.Range("M86").Value = -22180
$ws.Range("N86").Value = -9745.5

# Row 89
$ws.Range("H89").Value = 18787.715
$ws.Range("I89").Value = 23303
$ws.Range("J89").Value = 7499.5
$ws.Range("K89").Value = 116515
$ws.Range("L89").Value = 37497.5
$ws.Range("M89").Value = -110899
$ws.Range("N89").Value = -48729.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 22
$ws.Range("H22").Value = 7382.0625
$ws.Range("I22").Value = 9350.909
$ws.Range("K22").Value = 9350.909
$ws.Range("M22").Value = -9000.909

# Row 31
$ws.Range("H31").Value = 5955574
$ws.Range("I31").Value = 2088.5151
$ws.Range("J31").Value = 27785020
$ws.Range("K31").Value = 2088.5151
$ws.Range("L31").Value = 27785020
$ws.Range("M31").Value = -1793.5151
$ws.Range("N31").Value = -27785610

# Row 34
$ws.Range("H34").Value = 5955574
$ws.Range("I34").Value = 2088.5151
$ws.Range("J34").Value = 27785020
$ws.Range("K34").Value = 2088.5151
$ws.Range("L34").Value = 27785020
$ws.Range("M34").Value = -1886.5151
$ws.Range("N34").Value = -27785424

# Row 58
$ws.Range("H58").Value = 1157.2106
$ws.Range("I58").Value = 1183
$ws.Range("J58").Value = 1085
$ws.Range("K58").Value = 1183
$ws.Range("L58").Value = 1085
$ws.Range("M58").Value = -980
$ws.Range("N58").Value = -1491

# Row 105
$ws.Range("H105").Value = 17717.5
$ws.Range("I105").Value = 1258.8
$ws.Range("K105").Value = 1258.8
$ws.Range("M105").Value = 488.2

# Row 107
$ws.Range("H107").Value = 270
$ws.Range("I107").Value = 270
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 270
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1650
$ws.Range("N107").ClearContents()

# Row 122
$ws.Range("H122").Value = 1881379.5
$ws.Range("I122").Value = 1452.5264
$ws.Range("K122").Value = 4357.5792
$ws.Range("M122").Value = -1907.5792

# Row 132
$ws.Range("H132").Value = 185407.81
$ws.Range("I132").Value = 252807.38
$ws.Range("J132").Value = 5675.6665
$ws.Range("K132").Value = 758422.14
$ws.Range("L132").Value = 17026.9995
$ws.Range("M132").Value = -755892.14
$ws.Range("N132").Value = -22086.9995

# Row 134
$ws.Range("H134").Value = 2185.4583
$ws.Range("I134").Value = 1640.5714
$ws.Range("K134").Value = 4921.7142
$ws.Range("M134").Value = -2386.7142

# Row 136
$ws.Range("H136").Value = 1157.2106
$ws.Range("I136").Value = 1183
$ws.Range("J136").Value = 1085
$ws.Range("K136").Value = 3549
$ws.Range("L136").Value = 3255
$ws.Range("M136").Value = -999
$ws.Range("N136").Value = -8355

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 68
$ws.Range("H68").Value = 948.25
$ws.Range("J68").Value = 897
$ws.Range("L68").Value = 2691
$ws.Range("N68").Value = -4313

# Row 71
$ws.Range("H71").Value = 948.25
$ws.Range("J71").Value = 897
$ws.Range("L71").Value = 8073
$ws.Range("N71").Value = -16185

# Row 121
$ws.Range("H121").Value = 955.44446
$ws.Range("I121").Value = 1139.8
$ws.Range("K121").Value = 3419.4
$ws.Range("M121").Value = -2109.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 33
$ws.Range("H33").Value = 24999.5
$ws.Range("I33").Value = 24999.5
$ws.Range("K33").Value = 24999.5
$ws.Range("M33").Value = -24747.5

# Row 122
$ws.Range("H122").Value = 55557492
$ws.Range("I122").Value = 1634.8572
$ws.Range("K122").Value = 4904.571599999999
$ws.Range("M122").Value = -2454.571599999999

# Row 126
$ws.Range("H126").Value = 2549.4285
$ws.Range("I126").Value = 2474.3333
$ws.Range("K126").Value = 7422.999899999999
$ws.Range("M126").Value = -4952.999899999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 40
$ws.Range("H40").Value = 4674.5312
$ws.Range("I40").Value = 4371.8076
$ws.Range("K40").Value = 4371.8076
$ws.Range("M40").Value = -4235.8076

# Row 41
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

# Row 122
$ws.Range("H122").Value = 2979929.5
$ws.Range("J122").Value = 11367936
$ws.Range("L122").Value = 34103808
$ws.Range("N122").Value = -34108708

# Row 132
$ws.Range("H132").Value = 4912.857
$ws.Range("I132").Value = 2303
$ws.Range("J132").Value = 10607.091
$ws.Range("K132").Value = 6909
$ws.Range("L132").Value = 31821.273
$ws.Range("M132").Value = -4379
$ws.Range("N132").Value = -36881.273

# Row 136
$ws.Range("H136").Value = 10616.125
$ws.Range("I136").Value = 4494.75
$ws.Range("K136").Value = 13484.25
$ws.Range("M136").Value = -10934.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 122
$ws.Range("H122").Value = 11768196
$ws.Range("I122").Value = 1727.3636
$ws.Range("J122").Value = 33340056
$ws.Range("K122").Value = 5182.0908
$ws.Range("L122").Value = 100020168
$ws.Range("M122").Value = -2732.0908
$ws.Range("N122").Value = -100025068
